$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 48, shifting existing rows 48-78 down to 49-79
$ws.Rows.Item(48).Insert()

# Populate the newly inserted row 48 with the new weekly data entry
$ws.Range("A48").Value = 7
$ws.Range("B48").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C48").Value = "Ñuble"
$ws.Range("D48").Value = 44846
$ws.Range("E48").Value = 16
$ws.Range("F48").Value = 100112022
$ws.Range("G48").Value = "Arveja Verde"
$ws.Range("H48").Value = "Perfection"
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value = 60
$ws.Range("K48").Value = 25000
$ws.Range("L48").Value = 26000
$ws.Range("M48").Value = 25500
$ws.Range("N48").Value = "$/malla 25 kilos"
$ws.Range("O48").Value = "Provincia de Limarí"
$ws.Range("P48").Value = 1020
$ws.Range("Q48").Value = 25
$ws.Range("R48").Value = "Hortaliza"
